$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6753301551942219
$ws.Range("C2").Value = 10.29869402782916
$ws.Range("D2").Value = 337.1190423067083
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("G2").Value = 993.4203433196917

$ws.Range("B3").Value = 3.230985683306322
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 26.21740644021617
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 39.7764191927396
